# Update cryptos list values per upstream diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.890.98"
$ws.Range("E2").Value = "'  +2.72%  "
$ws.Range("D3").Value = "'1.668.28"
$ws.Range("E3").Value = "'  -0.46%  "
$ws.Range("E4").Value = "'  -0.08%  "
$ws.Range("D5").Value = "'214.75"
$ws.Range("E5").Value = "'  +0.17%  "
$ws.Range("E6").Value = "'  -0.29%  "
$ws.Range("E7").Value = "'  -0.07%  "
$ws.Range("D8").Value = "'23.64"
$ws.Range("E8").Value = "'  +3.52%  "
$ws.Range("E9").Value = "'  -0.10%  "
$ws.Range("E10").Value = "'  +0.03%  "
$ws.Range("D11").Value = "'0.0879"
$ws.Range("E11").Value = "'  -1.23%  "
$ws.Range("D12").Value = "'1.904.77"
$ws.Range("E12").Value = "'  -0.51%  "
$ws.Range("D13").Value = "'1.669.09"
$ws.Range("E13").Value = "'  -1.05%  "
$ws.Range("E14").Value = "'  -1.34%  "
$ws.Range("E15").Value = "'  -1.48%  "
$ws.Range("D16").Value = "'66.05"
$ws.Range("E16").Value = "'  -0.74%  "
$ws.Range("D17").Value = "'251.82"
$ws.Range("E17").Value = "'  +7.07%  "
$ws.Range("D18").Value = "'27.864.98"
$ws.Range("E18").Value = "'  +2.74%  "
$ws.Range("D19").Value = "'0.0₃0732"
$ws.Range("E19").Value = "'  -1.14%  "
$ws.Range("D20").Value = "'7.57"
$ws.Range("E20").Value = "'  -4.09%  "
$ws.Range("D21").Value = "'1.00"
$ws.Range("E21").Value = "'  -0.14%  "
$ws.Range("D22").Value = "'4.48"
$ws.Range("E22").Value = "'  -1.52%  "
$ws.Range("D23").Value = "'9.35"
$ws.Range("E23").Value = "'  -2.02%  "
$ws.Range("E24").Value = "'  -1.68%  "
$ws.Range("D25").Value = "'146.69"
$ws.Range("E25").Value = "'  -1.05%  "
$ws.Range("D26").Value = "'7.23"
$ws.Range("E26").Value = "'  -3.20%  "
$ws.Range("D27").Value = "'16.30"
$ws.Range("E27").Value = "'  -0.57%  "
$ws.Range("E28").Value = "'  -0.34%  "
$ws.Range("E29").Value = "'  -0.04%  "
$ws.Range("E30").Value = "'  +5.90%  "
$ws.Range("D31").Value = "'0.0502"
$ws.Range("E31").Value = "'  +0.26%  "
$ws.Range("E32").Value = "'  -0.41%  "
$ws.Range("E33").Value = "'  -2.61%  "
$ws.Range("D34").Value = "'1.427.30"
$ws.Range("E34").Value = "'  -7.24%  "
$ws.Range("E35").Value = "'  -5.79%  "
$ws.Range("E36").Value = "'  +0.05%  "
$ws.Range("D37").Value = "'0.930"
$ws.Range("E37").Value = "'  -1.51%  "
$ws.Range("E38").Value = "'  -4.13%  "
$ws.Range("E39").Value = "'  -1.16%  "
$ws.Range("B40").Value = "Aave"
$ws.Range("C40").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D40").Value = "'69.58"
$ws.Range("E40").Value = "'  -0.25%  "
$ws.Range("B41").Value = "WEMIXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D41").Value = "'1.03"
$ws.Range("E41").Value = "'  -3.06%  "
$ws.Range("E42").Value = "'  +0.02%  "
$ws.Range("D43").Value = "'2.22"
$ws.Range("E43").Value = "'  -1.44%  "
$ws.Range("D44").Value = "'1.812.87"
$ws.Range("E44").Value = "'  -0.53%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").Value = "'5.39"
$ws.Range("E45").Value = "'  -6.70%  "
$ws.Range("B46").Value = "TrustWalletToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D46").Value = "'0.793"
$ws.Range("E46").Value = "'  +1.40%  "
$ws.Range("D47").Value = "'1.72"
$ws.Range("E47").Value = "'  +5.15%  "
$ws.Range("D48").Value = "'88.94"
$ws.Range("E48").Value = "'  -0.83%  "
$ws.Range("E49").Value = "'  -0.55%  "
$ws.Range("E50").Value = "'  -2.00%  "
$ws.Range("D51").Value = "'7.81"
$ws.Range("E51").Value = "'  -4.83%  "
